# Update "Estado de Cuenta" table (Hoja1) with refreshed period range and
# salary values: the mora table is re-sorted ascending by period (1611..2003)
# and the "Salario Basico" column is refreshed to the new base salary.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

$periods = @(
    "1611","1612",
    "1701","1702","1703","1704","1705","1706","1707","1708","1709","1710","1711","1712",
    "1801","1802","1803","1804","1805","1806","1807","1808","1809","1810","1811","1812",
    "1901","1902","1903","1904","1905","1906","1907","1908","1909","1910","1911","1912",
    "2001","2002","2003"
)

$newSalarioBasico = 781242

$startRow = 16
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = $startRow + $i
    $period = $periods[$i]

    if ([int]$period -le 1808) {
        $valorMora = 27580
    } else {
        $valorMora = 31249
    }

    $ws.Range("E$row").Value = $period
    $ws.Range("F$row").Value = $valorMora
    $ws.Range("G$row").Value = $newSalarioBasico
}
